# Update existing averages with the latest computed values (more games added)

$wb = $excel.ActiveWorkbook

# Sheet: Narrator Votes Averages
$ws1 = $wb.Worksheets.Item("Narrator Votes Averages")
$ws1.Range("B2").Value = 53.33333333333333
$ws1.Range("C2").Value = 9.999999999999998
$ws1.Range("C3").Value = 50.41666666666666
$ws1.Range("B4").Value = 13.33333333333333
$ws1.Range("C4").Value = 39.58333333333333

# Sheet: Votes Not Narrator Averages
$ws2 = $wb.Worksheets.Item("Votes Not Narrator Averages")
$ws2.Range("B2").Value = 15.37995337995338
$ws2.Range("C2").Value = 44.1899766899767

# Sheet: Correct Votes Averages
$ws3 = $wb.Worksheets.Item("Correct Votes Averages")
$ws3.Range("B2").Value = 27.44755244755245
$ws3.Range("C2").Value = 31.79845519977098

# Add a new sheet with winners statistics, placed after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Winners Statistics"

$ws4.Range("A1").Value = "Player"
$ws4.Range("B1").Value = "Winner Percent"
$ws4.Range("A2").Value = "GPT"
$ws4.Range("B2").Value = 0
$ws4.Range("A3").Value = "Bot"
$ws4.Range("B3").Value = 100

# Match the header formatting (bold, bordered, centered) used on the other sheets
$ws1.Range("B1").Copy()
$ws4.Range("A1:B1").PasteSpecial(-4122)
